# Apply the "financial-system" StructureDefinition metadata refresh:
#  - URL moves from ibm.com to linuxforhealth.org
#  - Version bump 7.0.0 -> 8.0.0
#  - Date bump
#  - Publisher renamed from "Alvearie Team" to "LinuxForHealth Team"
#  - The root "Extension" row's rolled-up Constraint(s) text is cleared
#    (it now only shows up on the Extension.extension child row)

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/financial-system"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element; its "Constraint(s)" column (AI)
# no longer repeats the ele-1/ext-1 text (that now lives solely on the
# Extension.extension child row further down).
$elements.Range("AI2").Value = ""

# Row 5 (Extension.url) carries the same canonical URL as its "Fixed
# Value" (column Q); keep it in sync with the Metadata sheet's URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/financial-system"
